$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gains one new weekly record. It is inserted as row 306, which
# pushes the former rows 306-329 down to 307-330 (dimension grows from
# A1:R329 to A1:R330).
$ws.Rows(306).Insert()

# The newly-inserted row 306 shares most attributes (Mercado ID, Mercado,
# Región, Codreg, Categoría ID, Categoría, Variedad, Unidad de
# comercialización, Kg o Unidades, Clasificación) with the record that was
# originally in row 306 (now shifted to row 307) - copy them across.
$ws.Range("A306").Value = $ws.Range("A307").Value2
$ws.Range("B306").Value = $ws.Range("B307").Text
$ws.Range("C306").Value = $ws.Range("C307").Text
$ws.Range("E306").Value = $ws.Range("E307").Value2
$ws.Range("F306").Value = $ws.Range("F307").Value2
$ws.Range("G306").Value = $ws.Range("G307").Text
$ws.Range("H306").Value = $ws.Range("H307").Text
$ws.Range("N306").Value = $ws.Range("N307").Text
$ws.Range("Q306").Value = $ws.Range("Q307").Value2
$ws.Range("R306").Value = $ws.Range("R307").Text

# New weekly values for the new record.
$ws.Range("D306").Value = 45106
$ws.Range("I306").Value = "1a (guarda)"
$ws.Range("J306").Value = 250
$ws.Range("K306").Value = 400
$ws.Range("L306").Value = 450
$ws.Range("M306").Value = 430
$ws.Range("O306").Value = "Región del Maule"
$ws.Range("P306").Value = 430
